$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.606.99"
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = "'1.927.26"
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = '  +0.79%  '
$ws.Range("D5").Value = "'326.48"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = "'1.012"
$ws.Range("D7").Value = "'0.4826"
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").Value = "'0.4059"
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").Value = "'0.08195"
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("E10").Value = '  -0.63%  '
$ws.Range("D11").Value = "'23.75"
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "'1.915.41"
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'6.069"
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'7.289"
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").Value = "'91.49"
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = "'0.06858"
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").Value = "'0.00001039"
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = "'17.62"
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").Value = "'29.614.52"
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = "'5.656"
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").Value = "'11.94"
$ws.Range("E23").Value = '  +1.68%  '
$ws.Range("D24").Value = "'2.196"
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("D25").Value = "'2.115.56"
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("D26").Value = "'156.34"
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = "'6.362"
$ws.Range("E27").Value = '  -3.09%  '
$ws.Range("D28").Value = "'20.03"
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("E29").Value = '  -1.83%  '
$ws.Range("D30").Value = "'120.76"
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = "'1.004"
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("D32").Value = "'0.09605"
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").Value = "'5.617"
$ws.Range("E33").Value = '  +1.63%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = "'1.390"
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("D36").Value = "'0.06520"
$ws.Range("E36").Value = '  +6.29%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = "'1.219"
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("D39").Value = "'0.5939"
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("D40").Value = "'10.72"
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("D41").Value = "'7.863"
$ws.Range("E41").Value = '  -2.42%  '
$ws.Range("D42").Value = "'0.1845"
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").Value = "'2.490"
$ws.Range("E43").Value = '  +3.38%  '
$ws.Range("D44").Value = "'1.244"
$ws.Range("E44").Value = '  -2.73%  '
$ws.Range("D45").Value = "'0.07549"
$ws.Range("E46").Value = '  -1.46%  '
$ws.Range("D47").Value = "'0.5553"
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").Value = "'1.958"
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").Value = "'118.30"
$ws.Range("E49").Value = '  +1.29%  '
$ws.Range("D50").Value = "'2.432"
$ws.Range("E50").Value = '  +0.70%  '
$ws.Range("D51").Value = "'71.97"
$ws.Range("E51").Value = '  -1.16%  '
